$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary table values ---
$ws.Range("D2").Value = 2
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = 0

$ws.Range("D3").Value = 1
$ws.Range("J3").Value = $null
$ws.Range("M3").Value = $null

$ws.Range("G4").Value = 2
$ws.Range("J4").Value = $null
$ws.Range("M4").Value = $null

# --- Rewrite the Room Coordinates list (column A) ---
# The write order below matters: it determines the order new strings are
# interned into the shared string table, matching the target workbook.
$ws.Cells.Item(7, 1).Value = "(3,5)"
$ws.Cells.Item(11, 1).Value = "(3,3)"
$ws.Cells.Item(12, 1).Value = "(2,3)"
$ws.Cells.Item(13, 1).Value = "(2,2)"
$ws.Cells.Item(14, 1).Value = "(3,2)"
$ws.Cells.Item(3, 1).Value = "(0,2)"
$ws.Cells.Item(4, 1).Value = "(1,2)"
$ws.Cells.Item(5, 1).Value = "(1,4)"
$ws.Cells.Item(6, 1).Value = "(2,4)"
$ws.Cells.Item(8, 1).Value = "(3,4)"
$ws.Cells.Item(10, 1).Value = "(4,3)"
$ws.Cells.Item(17, 1).Value = "(4,1)"
$ws.Cells.Item(18, 1).Value = "(1,1)"
$ws.Cells.Item(19, 1).Value = "(1,0)"

# Rows re-using strings already introduced above
$ws.Cells.Item(9, 1).Value = "(4,4)"
$ws.Cells.Item(15, 1).Value = "(3,3)"
$ws.Cells.Item(16, 1).Value = "(4,3)"

# --- Selection ---
$ws.Range("G5").Select()
